$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "19.00", "2.538.34") are preserved exactly, then clear the
    # temporary number-format override so the cell keeps its original style.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "64.918.58"
$ws.Range("E2").Value = "  +3.40%  "
Set-TextValue "D3" "2.534.86"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "581.57"
$ws.Range("E5").Value = "  +1.16%  "
Set-TextValue "D6" "153.43"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.18%  "
Set-TextValue "D9" "2.538.34"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  -0.15%  "
Set-TextValue "D14" "29.34"
$ws.Range("E14").Value = "  +0.46%  "
Set-TextValue "D15" "0.0000182"
$ws.Range("E15").Value = "  +2.40%  "
Set-TextValue "D16" "2.995.02"
$ws.Range("E16").Value = "  +2.72%  "
Set-TextValue "D17" "64.776.52"
$ws.Range("E17").Value = "  +3.30%  "
Set-TextValue "D18" "2.538.04"
$ws.Range("E18").Value = "  +2.99%  "
Set-TextValue "D19" "7.99"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("E21").Value = "  +3.76%  "
Set-TextValue "D22" "329.87"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +1.05%  "
Set-TextValue "D26" "65.90"
$ws.Range("E26").Value = "  +0.54%  "
Set-TextValue "D27" "649.59"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  +7.29%  "
$ws.Range("E30").Value = "  +4.92%  "
$ws.Range("E31").Value = "  +0.78%  "
Set-TextValue "D32" "8.09"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +1.24%  "
Set-TextValue "D37" "4.87"
$ws.Range("E37").Value = "  +2.73%  "
$ws.Range("E38").Value = "  +4.69%  "
Set-TextValue "D39" "155.17"
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("E40").Value = "  +2.23%  "
Set-TextValue "D41" "19.00"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  +5.21%  "
Set-TextValue "D44" "163.03"
$ws.Range("E44").Value = "  +6.42%  "
$ws.Range("E45").Value = "  -0.01%  "
Set-TextValue "D46" "0.0₆0305"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("E48").Value = "  +2.29%  "
Set-TextValue "D49" "21.67"
$ws.Range("E49").Value = "  +5.99%  "
Set-TextValue "D50" "0.624"
$ws.Range("E50").Value = "  +2.57%  "
Set-TextValue "D51" "0.0519"
$ws.Range("E51").Value = "  +1.58%  "
